$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Shen, Jiayi" requested no faculty members at all -- remove her row entirely
# (was row 21: Shen, Jiayi / 0 / 1 / Emerson, Minin, Weining Shen, Swarup, Thornton, Zhaoxia Yu)
$ws.Rows.Item(21).Delete()

# After the deletion the remaining rows shift up by one. Two more students turn
# out to have requested no faculty -- clear their faculty-name column instead
# of deleting their rows (they still have B/C data).
$ws.Range("D7").ClearContents()   # Du, Mingyu -- was "Allard"
$ws.Range("D25").ClearContents()  # Xiang, Yankai (Mark) -- was "Allard"

# A couple of faculty lists gained an extra name
$ws.Range("D11").Value = "Allard, Collins, Komarova, Mjolsness, Siryaporn, Stern , Lander"
$ws.Range("D13").Value = "Cho, Cinquin, Lander, Gardiner, Schilling, Nie, Arora, Downing"

# Xu, Angela is now flagged with an asterisk
$ws.Range("B26").Value = 1

# Cosmetic: widen the faculty-names column so long lists are readable, and
# restore the previous zoom/selection state
$ws.Columns.Item(4).ColumnWidth = 97
$excel.ActiveWindow.Zoom = 91
$ws.Range("D14").Select()

# Keep the sheet's recorded sort-state range in sync with the now-smaller
# data range (one row shorter after the deletion above)
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A1"))
$sortObj.SetRange($ws.Range("A1:G27"))
$sortObj.Header = 1
$sortObj.Apply()
